$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 222, shifting existing rows 222:325 down to 223:326.
$ws.Rows(222).Insert()

# Populate the newly inserted row 222 with the new data record.
$ws.Cells.Item(222, 1).Value  = 4
$ws.Cells.Item(222, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(222, 3).Value  = "Los Lagos"
$ws.Cells.Item(222, 4).Value  = 44875
$ws.Cells.Item(222, 5).Value  = 10
$ws.Cells.Item(222, 6).Value  = "Fruta"
$ws.Cells.Item(222, 7).Value  = 100108
$ws.Cells.Item(222, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(222, 9).Value  = 100108005
$ws.Cells.Item(222, 10).Value = "Piña"
$ws.Cells.Item(222, 11).Value = "Caramelo"
$ws.Cells.Item(222, 12).Value = "Segunda"
$ws.Cells.Item(222, 13).Value = 200
$ws.Cells.Item(222, 14).Value = 34000
$ws.Cells.Item(222, 15).Value = 35000
$ws.Cells.Item(222, 16).Value = 34500
$ws.Cells.Item(222, 17).Value = "$/caja 14 unidades"
$ws.Cells.Item(222, 18).Value = "Ecuador"
$ws.Cells.Item(222, 19).Value = 2464
$ws.Cells.Item(222, 20).Value = 14

# Make sure the date cell keeps the same date style used by the rest of column D.
$ws.Cells.Item(222, 4).NumberFormat = $ws.Cells.Item(221, 4).NumberFormat
